$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C14: "grupa in [110…938]" -> "grupa in [111…938]"
$ws.Range("C14").Value = "grupa in [111…938]"

# Update D15: "grupa <110" -> "grupa <111"
$ws.Range("D15").Value = "grupa <111"

# Update D16: "grupa >938" -> "grupa >937"
$ws.Range("D16").Value = "grupa >937"

# A16 gets a new value 8
$ws.Range("A16").Value = 8

# Update H17:H22 boundary values
$ws.Range("H17").Value = 111
$ws.Range("H18").Value = 110
$ws.Range("H19").Value = 112
$ws.Range("H20").Value = 937
$ws.Range("H21").Value = 936
$ws.Range("H22").Value = 938

# Update selection/view
$ws.Range("A6").Select()
